$wb = $excel.ActiveWorkbook

# Overview sheet: update the d41011ff... row (row 3) status + handoff date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 14:51:02"

# zh-cn sheet: update the d41011ff... row (row 3) status + handoff datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-22 14:50:58"

# de-de sheet: update the d41011ff... row (row 3) status + handoff datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-22 14:51:02"
